$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.57873595852908
$ws.Range("C2").Value = 0.3354286270085822
$ws.Range("E2").Value = 0.151771830841156
$ws.Range("F2").Value = 4.253010948229331
$ws.Range("G2").Value = 0.002490816268178767
$ws.Range("J2").Value = 0.2876175167626656
$ws.Range("B3").Value = 1.4637037638214
$ws.Range("C3").Value = 0.303023540219499
$ws.Range("E3").Value = 0.1566965069044493
$ws.Range("F3").Value = 4.072601711761536
$ws.Range("G3").Value = 0.002497859305758434
$ws.Range("J3").Value = 0.2676841926982121
$ws.Range("B4").Value = 1.394403658732756
$ws.Range("C4").Value = 0.2834058188777817
$ws.Range("E4").Value = 0.1599641161839074
$ws.Range("F4").Value = 3.963759742175597
$ws.Range("G4").Value = 0.002502401404681309
$ws.Range("J4").Value = 0.2555252067487288
$ws.Range("B5").Value = 1.366493758243053
$ws.Range("C5").Value = 0.2754803192025008
$ws.Range("E5").Value = 0.1613566576718988
$ws.Range("F5").Value = 3.919879141737908
$ws.Range("G5").Value = 0.00250430730432379
$ws.Range("J5").Value = 0.2505889248810291
$ws.Range("B6").Value = 1.361879183926533
$ws.Range("C6").Value = 0.2741684186519819
$ws.Range("E6").Value = 0.161591560111864
$ws.Range("F6").Value = 3.912621046881412
$ws.Range("G6").Value = 0.002504627103344691
$ws.Range("J6").Value = 0.2497703375623246
$ws.Range("B7").Value = 1.394025923277411
$ws.Range("C7").Value = 0.2832986556259414
$ws.Range("E7").Value = 0.159982650098982
$ws.Range("F7").Value = 3.963166051985581
$ws.Range("G7").Value = 0.002502426885561594
$ws.Range("J7").Value = 0.255458561041749
$ws.Range("B8").Value = 1.538794301553821
$ws.Range("C8").Value = 0.3241966108668635
$ws.Range("E8").Value = 0.153419036962628
$ws.Range("F8").Value = 4.190397979588312
$ws.Range("G8").Value = 0.002493199673128926
$ws.Range("J8").Value = 0.2807269267806447
$ws.Range("B9").Value = 1.833440836133946
$ws.Range("C9").Value = 0.4066769682765425
$ws.Range("E9").Value = 0.1424986742728045
$ws.Range("F9").Value = 4.651880069963624
$ws.Range("G9").Value = 0.002476821595405415
$ws.Range("J9").Value = 0.330986641918372
$ws.Range("B10").Value = 2.056785614688124
$ws.Range("C10").Value = 0.4687631005024855
$ws.Range("E10").Value = 0.135687563439074
$ws.Range("F10").Value = 5.001439641024945
$ws.Range("G10").Value = 0.00246582029711114
$ws.Range("J10").Value = 0.3684485245056521
$ws.Range("B11").Value = 2.159953252807441
$ws.Range("C11").Value = 0.4973530836620625
$ws.Range("E11").Value = 0.1328573576333909
$ws.Range("F11").Value = 5.162929580383945
$ws.Range("G11").Value = 0.002461036366824866
$ws.Range("J11").Value = 0.3856303133860024
$ws.Range("B12").Value = 2.199250964670512
$ws.Range("C12").Value = 0.5082309693091247
$ws.Range("E12").Value = 0.1318246491685393
$ws.Range("F12").Value = 5.224451982200833
$ws.Range("G12").Value = 0.002459256294719139
$ws.Range("J12").Value = 0.3921586130541073
$ws.Range("B13").Value = 2.190777171755826
$ws.Range("C13").Value = 0.5058859060128498
$ws.Range("E13").Value = 0.1320453180642396
$ws.Range("F13").Value = 5.211185365516997
$ws.Range("G13").Value = 0.002459638267683823
$ws.Range("J13").Value = 0.3907516231178647
$ws.Range("B14").Value = 2.163181652901869
$ws.Range("C14").Value = 0.4982469726407999
$ws.Range("E14").Value = 0.1327716111858948
$ws.Range("F14").Value = 5.167983575652329
$ws.Range("G14").Value = 0.00246088928930326
$ws.Range("J14").Value = 0.386166950789061
$ws.Range("B15").Value = 2.146308753543678
$ws.Range("C15").Value = 0.4935746586553478
$ws.Range("E15").Value = 0.1332215833994219
$ws.Range("F15").Value = 5.141569793169367
$ws.Range("G15").Value = 0.002461659671458791
$ws.Range("J15").Value = 0.3833616154739445
$ws.Range("B16").Value = 2.050075325409296
$ws.Range("C16").Value = 0.466901791995042
$ws.Range("E16").Value = 0.1358779600944544
$ws.Range("F16").Value = 4.990936855613626
$ws.Range("G16").Value = 0.002466137358033135
$ws.Range("J16").Value = 0.3673286258045607
$ws.Range("B17").Value = 1.991444108048711
$ws.Range("C17").Value = 0.4506287888792144
$ws.Range("E17").Value = 0.1375765780340537
$ws.Range("F17").Value = 4.899171812083154
$ws.Range("G17").Value = 0.002468940618267525
$ws.Range("J17").Value = 0.3575300620856865
$ws.Range("B18").Value = 1.95786816299551
$ws.Range("C18").Value = 0.4413015197220034
$ws.Range("E18").Value = 0.1385787941106926
$ws.Range("F18").Value = 4.846622954160694
$ws.Range("G18").Value = 0.002470573758135278
$ws.Range("J18").Value = 0.3519072751802241
$ws.Range("B19").Value = 1.946525052434652
$ws.Range("C19").Value = 0.4381490053089578
$ws.Range("E19").Value = 0.1389224442939057
$ws.Range("F19").Value = 4.828870233386766
$ws.Range("G19").Value = 0.002471130286987737
$ws.Range("J19").Value = 0.3500056905572677
$ws.Range("B20").Value = 1.997670232640587
$ws.Range("C20").Value = 0.4523576993031497
$ws.Range("E20").Value = 0.1373931442027754
$ws.Range("F20").Value = 4.908916249507229
$ws.Range("G20").Value = 0.002468640057435412
$ws.Range("J20").Value = 0.3585717686587202
$ws.Range("B21").Value = 2.171280832693867
$ws.Range("C21").Value = 0.5004893021941825
$ws.Range("E21").Value = 0.1325572182949486
$ws.Range("F21").Value = 5.180662843092591
$ws.Range("G21").Value = 0.002460520980872771
$ws.Range("J21").Value = 0.3875129707836606
$ws.Range("B22").Value = 2.286090600868647
$ws.Range("C22").Value = 0.5322469400640557
$ws.Range("E22").Value = 0.1296243778688471
$ws.Range("F22").Value = 5.360426346759709
$ws.Range("G22").Value = 0.002455398185681445
$ws.Range("J22").Value = 0.4065563822430533
$ws.Range("B23").Value = 2.224689674028355
$ws.Range("C23").Value = 0.5152692125128624
$ws.Range("E23").Value = 0.1311686972284765
$ws.Range("F23").Value = 5.264280719889143
$ws.Range("G23").Value = 0.002458115603369706
$ws.Range("J23").Value = 0.3963801898119641
$ws.Range("B24").Value = 1.994854991181228
$ws.Range("C24").Value = 0.4515759708842779
$ws.Range("E24").Value = 0.1374759947538475
$ws.Range("F24").Value = 4.904510142253798
$ws.Range("G24").Value = 0.002468775874018332
$ws.Range("J24").Value = 0.358100780427236
$ws.Range("B25").Value = 1.752548033529877
$ws.Range("C25").Value = 0.3841098247823993
$ws.Range("E25").Value = 0.145241762388034
$ws.Range("F25").Value = 4.525252178146701
$ws.Range("G25").Value = 0.002481070061226843
$ws.Range("J25").Value = 0.3173035326923639
